$wb = $excel.ActiveWorkbook

# Data for the new row (row 78) to append to each of the 4 worksheets.
# Each entry: SheetName, B, C, D, E, F, G, H, I
$rows = @(
    @{ Sheet = "FE_LFT_#1"; B = "0x01,0x7c"; C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"; D = "0x01,0x24"; E = "0xf"; F = 380; G = "7.598631275147109e+23"; H = 292; I = 15 },
    @{ Sheet = "FE_LFT_#2"; B = "0x01,0x90"; C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"; D = "0x01,0x30"; E = "0xe"; F = 400; G = "5.68432987514711e+23";  H = 304; I = 14 },
    @{ Sheet = "FE_PLT_#1"; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x00,0x60"; E = "0x3"; F = 110; G = "5.68631262647114e+23"; H = 96;  I = 3 },
    @{ Sheet = "FE_PLT_#2"; B = "0x00,0x6e"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x00,0x5F"; E = "0x3"; F = 110; G = "9.85046333984776e+23"; H = 95;  I = 3 }
)

$newRow = 78
$dateValue = 45864.49234953704

foreach ($entry in $rows) {
    $ws = $wb.Worksheets.Item($entry.Sheet)

    # Column A keeps the same date/time style as the row above it (row 77).
    $ws.Cells.Item($newRow, 1).Value = $dateValue
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $entry.B
    $ws.Cells.Item($newRow, 3).Value = $entry.C
    $ws.Cells.Item($newRow, 4).Value = $entry.D
    $ws.Cells.Item($newRow, 5).Value = $entry.E
    $ws.Cells.Item($newRow, 6).Value = $entry.F
    $ws.Cells.Item($newRow, 7).Value = [double]$entry.G
    $ws.Cells.Item($newRow, 8).Value = $entry.H
    $ws.Cells.Item($newRow, 9).Value = $entry.I
}
